$d = $word.ActiveDocument

$d.Content.Find.Execute("34×48=1632", $true, $false, $false, $false, $false, $true, 1, $false, "36×30=1080", 2) | Out-Null
$d.Content.Find.Execute("61×24=1464", $true, $false, $false, $false, $false, $true, 1, $false, "21×47=987", 2) | Out-Null
$d.Content.Find.Execute("52×29=1508", $true, $false, $false, $false, $false, $true, 1, $false, "14×15=210", 2) | Out-Null
$d.Content.Find.Execute("91×35=3185", $true, $false, $false, $false, $false, $true, 1, $false, "81×53=4293", 2) | Out-Null
$d.Content.Find.Execute("60×16=960", $true, $false, $false, $false, $false, $true, 1, $false, "49×38=1862", 2) | Out-Null
$d.Content.Find.Execute("89×17=1513", $true, $false, $false, $false, $false, $true, 1, $false, "99×60=5940", 2) | Out-Null
$d.Content.Find.Execute("51×72=3672", $true, $false, $false, $false, $false, $true, 1, $false, "88×12=1056", 2) | Out-Null
$d.Content.Find.Execute("19×82=1558", $true, $false, $false, $false, $false, $true, 1, $false, "28×57=1596", 2) | Out-Null
$d.Content.Find.Execute("36×16=576", $true, $false, $false, $false, $false, $true, 1, $false, "39×88=3432", 2) | Out-Null
$d.Content.Find.Execute("20×54=1080", $true, $false, $false, $false, $false, $true, 1, $false, "37×53=1961", 2) | Out-Null
$d.Content.Find.Execute("62×81=5022", $true, $false, $false, $false, $false, $true, 1, $false, "71×96=6816", 2) | Out-Null
$d.Content.Find.Execute("65×79=5135", $true, $false, $false, $false, $false, $true, 1, $false, "69×96=6624", 2) | Out-Null
$d.Content.Find.Execute("44×12=528", $true, $false, $false, $false, $false, $true, 1, $false, "63×22=1386", 2) | Out-Null
$d.Content.Find.Execute("87×55=4785", $true, $false, $false, $false, $false, $true, 1, $false, "71×93=6603", 2) | Out-Null
$d.Content.Find.Execute("91×71=6461", $true, $false, $false, $false, $false, $true, 1, $false, "89×36=3204", 2) | Out-Null
$d.Content.Find.Execute("29×82=2378", $true, $false, $false, $false, $false, $true, 1, $false, "39×76=2964", 2) | Out-Null
$d.Content.Find.Execute("12×88=1056", $true, $false, $false, $false, $false, $true, 1, $false, "70×70=4900", 2) | Out-Null
$d.Content.Find.Execute("79×87=6873", $true, $false, $false, $false, $false, $true, 1, $false, "18×82=1476", 2) | Out-Null
$d.Content.Find.Execute("63×18=1134", $true, $false, $false, $false, $false, $true, 1, $false, "36×86=3096", 2) | Out-Null
$d.Content.Find.Execute("84×81=6804", $true, $false, $false, $false, $false, $true, 1, $false, "83×78=6474", 2) | Out-Null
$d.Content.Find.Execute("14×18=252", $true, $false, $false, $false, $false, $true, 1, $false, "49×54=2646", 2) | Out-Null
$d.Content.Find.Execute("87×57=4959", $true, $false, $false, $false, $false, $true, 1, $false, "12×80=960", 2) | Out-Null
$d.Content.Find.Execute("85×42=3570", $true, $false, $false, $false, $false, $true, 1, $false, "76×31=2356", 2) | Out-Null
$d.Content.Find.Execute("74×88=6512", $true, $false, $false, $false, $false, $true, 1, $false, "66×22=1452", 2) | Out-Null
$d.Content.Find.Execute("41×46=1886", $true, $false, $false, $false, $false, $true, 1, $false, "80×67=5360", 2) | Out-Null
